{"js": "// \"ferdig oppgave 7 +\"\n// Append a 3x2 results table (\"Wins\"/\"Lose\"/\"Remi\" header row, with a\n// second row of numbers 704/296/1600) right after the existing body\n// content, followed by a trailing empty paragraph, and before the\n// section properties.\n\nconst tableXml =\n  '<w:tbl>' +\n    '<w:tblPr>' +\n      '<w:tblBorders>' +\n        '<w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n        '<w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n        '<w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n        '<w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n        '<w:insideH w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n        '<w:insideV w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n      '</w:tblBorders>' +\n      '<w:tblCellMar>' +\n        '<w:left w:w=\"108\" w:type=\"dxa\"/>' +\n        '<w:right w:w=\"108\" w:type=\"dxa\"/>' +\n      '</w:tblCellMar>' +\n    '</w:tblPr>' +\n    '<w:tblGrid>' +\n      '<w:gridCol w:w=\"246\"/>' +\n      '<w:gridCol w:w=\"246\"/>' +\n      '<w:gridCol w:w=\"246\"/>' +\n    '</w:tblGrid>' +\n    '<w:tr>' +\n      '<w:tblPrEx>' +\n        '<w:tblBorders>' +\n          '<w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n          '<w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n          '<w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n          '<w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n          '<w:insideH w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n          '<w:insideV w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n        '</w:tblBorders>' +\n        '<w:tblCellMar>' +\n          '<w:left w:w=\"108\" w:type=\"dxa\"/>' +\n          '<w:right w:w=\"108\" w:type=\"dxa\"/>' +\n        '</w:tblCellMar>' +\n      '</w:tblPrEx>' +\n      '<w:tc>' +\n        '<w:tcPr><w:vAlign w:val=\"center\"/></w:tcPr>' +\n        '<w:p><w:r><w:t xml:space=\"preserve\">Wins </w:t></w:r></w:p>' +\n      '</w:tc>' +\n      '<w:tc>' +\n        '<w:tcPr><w:vAlign w:val=\"center\"/></w:tcPr>' +\n        '<w:p><w:r><w:t>Lose</w:t></w:r></w:p>' +\n      '</w:tc>' +\n      '<w:tc>' +\n        '<w:tcPr><w:vAlign w:val=\"center\"/></w:tcPr>' +\n        '<w:p><w:r><w:t>Remi</w:t></w:r></w:p>' +\n      '</w:tc>' +\n    '</w:tr>' +\n    '<w:tr>' +\n      '<w:tblPrEx>' +\n        '<w:tblW w:w=\"5000\" w:type=\"pct\"/>' +\n        '<w:tblCellMar>' +\n          '<w:left w:w=\"108\" w:type=\"dxa\"/>' +\n          '<w:right w:w=\"108\" w:type=\"dxa\"/>' +\n        '</w:tblCellMar>' +\n      '</w:tblPrEx>' +\n      '<w:trPr><w:trHeight w:hRule=\"exact\" w:val=\"2000\"/></w:trPr>' +\n      '<w:tc>' +\n        '<w:tcPr><w:textDirection w:val=\"lrTb\"/><w:vAlign w:val=\"center\"/></w:tcPr>' +\n        '<w:p><w:r><w:t>704</w:t></w:r></w:p>' +\n        '<w:p/>' +\n      '</w:tc>' +\n      '<w:tc>' +\n        '<w:tcPr><w:textDirection w:val=\"lrTb\"/><w:vAlign w:val=\"center\"/></w:tcPr>' +\n        '<w:p><w:r><w:t>296</w:t></w:r></w:p>' +\n        '<w:p/>' +\n      '</w:tc>' +\n      '<w:tc>' +\n        '<w:tcPr><w:textDirection w:val=\"lrTb\"/><w:vAlign w:val=\"center\"/></w:tcPr>' +\n        '<w:p><w:r><w:t>1600</w:t></w:r></w:p>' +\n        '<w:p/>' +\n      '</w:tc>' +\n    '</w:tr>' +\n  '</w:tbl>' +\n  '<w:p/>';\n\nconst ooxmlPackage =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' + tableXml + '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\n// Anchor at the very end of the document body (just before sectPr) and\n// insert the table + trailing paragraph right after it.\nconst body = context.document.body;\nconst endRange = body.getRange(Word.RangeLocation.end);\nendRange.insertOoxml(ooxmlPackage, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# \"ferdig oppgave 7 +\"\n# Append a 3x2 results table (\"Wins\"/\"Lose\"/\"Remi\" header row, with a\n# second row of numbers 704/296/1600) right after the existing body\n# content, followed by a trailing empty paragraph, and before the\n# section properties.\n\n$d = $word.ActiveDocument\n\n# The raw WordprocessingML for the table plus the trailing empty\n# paragraph, exactly as it should land in word/document.xml.\n$tableXml = '<w:tbl>' +\n  '<w:tblPr>' +\n    '<w:tblBorders>' +\n      '<w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n      '<w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n      '<w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n      '<w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n      '<w:insideH w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n      '<w:insideV w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n    '</w:tblBorders>' +\n    '<w:tblCellMar>' +\n      '<w:left w:w=\"108\" w:type=\"dxa\"/>' +\n      '<w:right w:w=\"108\" w:type=\"dxa\"/>' +\n    '</w:tblCellMar>' +\n  '</w:tblPr>' +\n  '<w:tblGrid>' +\n    '<w:gridCol w:w=\"246\"/>' +\n    '<w:gridCol w:w=\"246\"/>' +\n    '<w:gridCol w:w=\"246\"/>' +\n  '</w:tblGrid>' +\n  '<w:tr>' +\n    '<w:tblPrEx>' +\n      '<w:tblBorders>' +\n        '<w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n        '<w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n        '<w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n        '<w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n        '<w:insideH w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n        '<w:insideV w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"auto\"/>' +\n      '</w:tblBorders>' +\n      '<w:tblCellMar>' +\n        '<w:left w:w=\"108\" w:type=\"dxa\"/>' +\n        '<w:right w:w=\"108\" w:type=\"dxa\"/>' +\n      '</w:tblCellMar>' +\n    '</w:tblPrEx>' +\n    '<w:tc>' +\n      '<w:tcPr><w:vAlign w:val=\"center\"/></w:tcPr>' +\n      '<w:p><w:r><w:t xml:space=\"preserve\">Wins </w:t></w:r></w:p>' +\n    '</w:tc>' +\n    '<w:tc>' +\n      '<w:tcPr><w:vAlign w:val=\"center\"/></w:tcPr>' +\n      '<w:p><w:r><w:t>Lose</w:t></w:r></w:p>' +\n    '</w:tc>' +\n    '<w:tc>' +\n      '<w:tcPr><w:vAlign w:val=\"center\"/></w:tcPr>' +\n      '<w:p><w:r><w:t>Remi</w:t></w:r></w:p>' +\n    '</w:tc>' +\n  '</w:tr>' +\n  '<w:tr>' +\n    '<w:tblPrEx>' +\n      '<w:tblW w:w=\"5000\" w:type=\"pct\"/>' +\n      '<w:tblCellMar>' +\n        '<w:left w:w=\"108\" w:type=\"dxa\"/>' +\n        '<w:right w:w=\"108\" w:type=\"dxa\"/>' +\n      '</w:tblCellMar>' +\n    '</w:tblPrEx>' +\n    '<w:trPr><w:trHeight w:hRule=\"exact\" w:val=\"2000\"/></w:trPr>' +\n    '<w:tc>' +\n      '<w:tcPr><w:textDirection w:val=\"lrTb\"/><w:vAlign w:val=\"center\"/></w:tcPr>' +\n      '<w:p><w:r><w:t>704</w:t></w:r></w:p>' +\n      '<w:p/>' +\n    '</w:tc>' +\n    '<w:tc>' +\n      '<w:tcPr><w:textDirection w:val=\"lrTb\"/><w:vAlign w:val=\"center\"/></w:tcPr>' +\n      '<w:p><w:r><w:t>296</w:t></w:r></w:p>' +\n      '<w:p/>' +\n    '</w:tc>' +\n    '<w:tc>' +\n      '<w:tcPr><w:textDirection w:val=\"lrTb\"/><w:vAlign w:val=\"center\"/></w:tcPr>' +\n      '<w:p><w:r><w:t>1600</w:t></w:r></w:p>' +\n      '<w:p/>' +\n    '</w:tc>' +\n  '</w:tr>' +\n'</w:tbl>' +\n'<w:p/>'\n\n$pkg = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $tableXml + '</w:body>' +\n      '</w:document>' +\n    '</pkg:xmlData>' +\n  '</pkg:part>' +\n'</pkg:package>'\n\n# Collapse to the very end of the document body (just before sectPr) and\n# insert the table + trailing paragraph there.\n$rng = $d.Content\n$rng.Collapse(0)\n$rng.InsertXML($pkg)\n"}
